$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPRINT BACKLOG 1")

# --- Row 11: "Nro de Horas" label + SUM formulas (burndown totals) ---
$ws.Range("B11").Value = "Nro de Horas"
$ws.Range("C11").Formula = "=SUM(C6:C10)"
$ws.Range("D11:I11").Formula = "=SUM(D6:D10)"

# --- New row 25: author addition ---
$ws.Range("A25").Value = "Carlos Zarate Carpio"

# --- Sprint Burndown line chart, built from the row 11 totals ---
$co = $ws.ChartObjects().Add(1276350, 2757487, 1343025, 2105025)
$chart = $co.Chart
$chart.ChartType = 4

$ser = $chart.SeriesCollection().NewSeries()
$ser.Name = "='SPRINT BACKLOG 1'!`$B`$11"
$ser.Values = $ws.Range("C11:I11")
$ser.XValues = $ws.Range("C4:I4")

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Sprint Burndown 1"

$chart.Axes(1).HasTitle = $true
$chart.Axes(1).AxisTitle.Text = "Nro Dias"

$chart.Axes(2).HasTitle = $true
$chart.Axes(2).AxisTitle.Text = "Nro de Horas"

$chart.HasLegend = $false

# --- Move selection the way the author's workbook shows it ---
$ws.Range("B21").Select()

Write-Host "edit applied"
